# Update Active_Outages.xlsx - 6/19/2025, 7:43:40 AM
# Refresh "Elapsed Duration(Hrs)" values on each region sheet and drop a
# resolved outage row (JED0123) from the R1 sheet.

$wb = $excel.ActiveWorkbook

# --- R1 --------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3944:57:47"
$ws1.Range("G3").Value = "84:30:25"
$ws1.Range("G4").Value = "107:30:25"
# Row 6 (JED0123, no open PCM ticket) has been resolved/removed.
$ws1.Rows.Item(6).Delete()

# --- R2 --------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12126:21:27"
$ws2.Range("G3").Value = "3256:04:56"
$ws2.Range("G4").Value = "494:16:30"

# --- R4 --------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2972:11:16"
$ws4.Range("G3").Value = "199:23:31"
$ws4.Range("G4").Value = "87:35:56"
$ws4.Range("G5").Value = "85:13:29"

# --- R5 --------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "446:10:15"

# --- R6 --------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "86:42:33"
